$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the mobile number values in column A (rows 2-5)
$ws.Range("A2").Value = 9363339066
$ws.Range("A3").Value = 9363339066
$ws.Range("A4").Value = 9363339066
$ws.Range("A5").Value = 9363339066

# Update the active selection to A6
$ws.Range("A6").Select()
